$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "V3.1"
$ws2 = $wb.Worksheets.Add($null, $ws)
$ws2.Name = "V3.2-evalboards"
Write-Host "index1=" $ws.Index
Write-Host "index2=" $ws2.Index
foreach ($s in $wb.Worksheets) {
  Write-Host $s.Name
}
